$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values per row (columns D, E, F, G, H)
$updates = @{
    4  = @(114368.1503112149, 0.002260488759696372, 0.186656093602205, -0.7238848619988547, 7.132001242965023)
    5  = @(115016.9714831865, -0.007010970750948553, 0.2217666072818137, -1.282494831271531, 10.34510309377354)
    7  = @(116336.8180687226, -0.01240109208217384, 0.1994258169348565, -0.2084160145127783, 4.205177023942746)
    9  = @(118310.5241376469, -0.05397337747656664, 0.306348362242314, -1.58321860560162, 10.08824214273261)
    10 = @(119711.8788105291, -0.1036865728320462, 0.4361538836730957, -1.874611108904448, 9.535701972788869)
    13 = @(112832.2978030683, 0.05976886482868678, 0.1393504333835252, -0.5297588083094721, 5.977643132738889)
    16 = @(112895.7243909183, 0.05224065531551625, 0.1481421723777032, -0.1969194871440355, 4.405082358658083)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}
